# Apply edits described by the diff:
# 1. Column G ("Fit") for rows 2-25 (Casacas) is cleared.
# 2. Column G for rows 26-37 (Jeans) changes text from "Pantalón Baggy" to "Baggy".
# 3. Column C width set to 43.1640625, column D width set to 22.5.
# 4. Selection changes to E34:E37.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the "Fit" values for rows 2-25 (Casacas/Poleras section)
$ws.Range("G2:G25").ClearContents()

# Update the "Fit" values for rows 26-37 (Jeans/Baggy section)
$ws.Range("G26:G37").Value = "Baggy"

# Set column widths
$ws.Columns.Item(3).ColumnWidth = 43.1640625
$ws.Columns.Item(4).ColumnWidth = 22.5

# Update the selection
$ws.Range("E34:E37").Select()
